$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the duplicate "Нетбэк 2020 прогноз" column (old column E) and the
# trailing empty column (old column N, which becomes column M after the
# first delete shifts everything left).
$ws.Columns("E").Delete()
$ws.Columns("M").Delete()

# Row heights: header row shrinks from 135 -> 60, sub-header row from 30 -> 15.
$ws.Rows(1).RowHeight = 60
$ws.Rows(2).RowHeight = 15

# A3: was a plain integer (8); now an actual date (01-Jan-2021), displayed
# with a short-date format.
$ws.Range("A3").NumberFormat = "mm-dd-yy"
$ws.Range("A3").Value = 44197

# F3 (old G3) was "=K3/J3" (a formula); replace with the new static figure
# the author typed in by hand.
$ws.Range("F3").Value = 65.6

# K3 (old L3) was "=G3*0.7" (a formula); replace with the new static figure
# the author typed in by hand.
$ws.Range("K3").Value = 45.9

# Restore the current selection to D1, matching the saved view state.
$ws.Range("D1").Select() | Out-Null
